# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (new "Feria Lagunitas de Puerto Montt" - Limón
# records dated 2021-09-21 / serial 44460) right before the existing row 173 data block,
# pushing the remainder of the table down by two rows and extending the sheet from
# A1:T261 to A1:T263.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 173:174 - shifts old rows 173..261 down to 175..263.
$ws.Rows("173:174").Insert()

# New row 173 - "1a amarillo"
$ws.Cells.Item(173, 1).Value = 4
$ws.Cells.Item(173, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(173, 3).Value = "Los Lagos"
$ws.Cells.Item(173, 4).Value = 44460
$ws.Cells.Item(173, 5).Value = 10
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100102
$ws.Cells.Item(173, 8).Value = "Cítricos"
$ws.Cells.Item(173, 9).Value = 100102003
$ws.Cells.Item(173, 10).Value = "Limón"
$ws.Cells.Item(173, 11).Value = "Sin especificar"
$ws.Cells.Item(173, 12).Value = "1a amarillo"
$ws.Cells.Item(173, 13).Value = 300
$ws.Cells.Item(173, 14).Value = 13000
$ws.Cells.Item(173, 15).Value = 13000
$ws.Cells.Item(173, 16).Value = 13000
$ws.Cells.Item(173, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(173, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(173, 19).Value = 722
$ws.Cells.Item(173, 20).Value = 18

# New row 174 - "2a amarillo"
$ws.Cells.Item(174, 1).Value = 4
$ws.Cells.Item(174, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(174, 3).Value = "Los Lagos"
$ws.Cells.Item(174, 4).Value = 44460
$ws.Cells.Item(174, 5).Value = 10
$ws.Cells.Item(174, 6).Value = "Fruta"
$ws.Cells.Item(174, 7).Value = 100102
$ws.Cells.Item(174, 8).Value = "Cítricos"
$ws.Cells.Item(174, 9).Value = 100102003
$ws.Cells.Item(174, 10).Value = "Limón"
$ws.Cells.Item(174, 11).Value = "Sin especificar"
$ws.Cells.Item(174, 12).Value = "2a amarillo"
$ws.Cells.Item(174, 13).Value = 300
$ws.Cells.Item(174, 14).Value = 9500
$ws.Cells.Item(174, 15).Value = 9500
$ws.Cells.Item(174, 16).Value = 9500
$ws.Cells.Item(174, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(174, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(174, 19).Value = 528
$ws.Cells.Item(174, 20).Value = 18

# Two new rows were appended at the very end (262, 263), duplicating the last
# observation pair (2021-06-02 / serial 44323, "1a plateado" + "2a plateado").
$ws.Cells.Item(262, 1).Value = 4
$ws.Cells.Item(262, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(262, 3).Value = "Los Lagos"
$ws.Cells.Item(262, 4).Value = 44323
$ws.Cells.Item(262, 5).Value = 10
$ws.Cells.Item(262, 6).Value = "Fruta"
$ws.Cells.Item(262, 7).Value = 100102
$ws.Cells.Item(262, 8).Value = "Cítricos"
$ws.Cells.Item(262, 9).Value = 100102003
$ws.Cells.Item(262, 10).Value = "Limón"
$ws.Cells.Item(262, 11).Value = "Sin especificar"
$ws.Cells.Item(262, 12).Value = "1a plateado"
$ws.Cells.Item(262, 13).Value = 1100
$ws.Cells.Item(262, 14).Value = 20000
$ws.Cells.Item(262, 15).Value = 21000
$ws.Cells.Item(262, 16).Value = 20500
$ws.Cells.Item(262, 17).Value = '$/malla 16 kilos'
$ws.Cells.Item(262, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(262, 19).Value = 1281
$ws.Cells.Item(262, 20).Value = 16

$ws.Cells.Item(263, 1).Value = 4
$ws.Cells.Item(263, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(263, 3).Value = "Los Lagos"
$ws.Cells.Item(263, 4).Value = 44323
$ws.Cells.Item(263, 5).Value = 10
$ws.Cells.Item(263, 6).Value = "Fruta"
$ws.Cells.Item(263, 7).Value = 100102
$ws.Cells.Item(263, 8).Value = "Cítricos"
$ws.Cells.Item(263, 9).Value = 100102003
$ws.Cells.Item(263, 10).Value = "Limón"
$ws.Cells.Item(263, 11).Value = "Sin especificar"
$ws.Cells.Item(263, 12).Value = "2a plateado"
$ws.Cells.Item(263, 13).Value = 400
$ws.Cells.Item(263, 14).Value = 17000
$ws.Cells.Item(263, 15).Value = 17000
$ws.Cells.Item(263, 16).Value = 17000
$ws.Cells.Item(263, 17).Value = '$/malla 16 kilos'
$ws.Cells.Item(263, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(263, 19).Value = 1062
$ws.Cells.Item(263, 20).Value = 16
